$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows at position 108, shifting existing rows 108-180 down to 117-189
$ws.Range("A108:I116").EntireRow.Insert()

# Force text format on B (date) and C (id) columns for the new rows so that
# date-like / numeric-like strings are not auto-converted by Excel
$textRng = $ws.Range("B108:C116")
$textRng.NumberFormat = "@"

# Row 108: 2019-11-18
$ws.Cells.Item(108,1).Value = 1574035200
$ws.Cells.Item(108,2).Value = "2019-11-18"
$ws.Cells.Item(108,3).Value = "0208"
$ws.Cells.Item(108,4).Value = "GREATEC"
$ws.Cells.Item(108,5).Value = 1.92
$ws.Cells.Item(108,6).Value = 1.96
$ws.Cells.Item(108,7).Value = 1.91
$ws.Cells.Item(108,8).Value = 1.95
$ws.Cells.Item(108,9).Value = 1363200

# Row 109: 2019-11-19
$ws.Cells.Item(109,1).Value = 1574121600
$ws.Cells.Item(109,2).Value = "2019-11-19"
$ws.Cells.Item(109,3).Value = "0208"
$ws.Cells.Item(109,4).Value = "GREATEC"
$ws.Cells.Item(109,5).Value = 1.96
$ws.Cells.Item(109,6).Value = 1.98
$ws.Cells.Item(109,7).Value = 1.92
$ws.Cells.Item(109,8).Value = 1.93
$ws.Cells.Item(109,9).Value = 1148700

# Row 110: 2019-11-20
$ws.Cells.Item(110,1).Value = 1574208000
$ws.Cells.Item(110,2).Value = "2019-11-20"
$ws.Cells.Item(110,3).Value = "0208"
$ws.Cells.Item(110,4).Value = "GREATEC"
$ws.Cells.Item(110,5).Value = 1.96
$ws.Cells.Item(110,6).Value = 1.96
$ws.Cells.Item(110,7).Value = 1.93
$ws.Cells.Item(110,8).Value = 1.93
$ws.Cells.Item(110,9).Value = 746900

# Row 111: 2019-11-21
$ws.Cells.Item(111,1).Value = 1574294400
$ws.Cells.Item(111,2).Value = "2019-11-21"
$ws.Cells.Item(111,3).Value = "0208"
$ws.Cells.Item(111,4).Value = "GREATEC"
$ws.Cells.Item(111,5).Value = 1.94
$ws.Cells.Item(111,6).Value = 1.94
$ws.Cells.Item(111,7).Value = 1.91
$ws.Cells.Item(111,8).Value = 1.92
$ws.Cells.Item(111,9).Value = 943000

# Row 112: 2019-11-22
$ws.Cells.Item(112,1).Value = 1574380800
$ws.Cells.Item(112,2).Value = "2019-11-22"
$ws.Cells.Item(112,3).Value = "0208"
$ws.Cells.Item(112,4).Value = "GREATEC"
$ws.Cells.Item(112,5).Value = 1.92
$ws.Cells.Item(112,6).Value = 1.95
$ws.Cells.Item(112,7).Value = 1.9
$ws.Cells.Item(112,8).Value = 1.92
$ws.Cells.Item(112,9).Value = 1405200

# Row 113: 2019-11-25
$ws.Cells.Item(113,1).Value = 1574640000
$ws.Cells.Item(113,2).Value = "2019-11-25"
$ws.Cells.Item(113,3).Value = "0208"
$ws.Cells.Item(113,4).Value = "GREATEC"
$ws.Cells.Item(113,5).Value = 1.92
$ws.Cells.Item(113,6).Value = 1.93
$ws.Cells.Item(113,7).Value = 1.89
$ws.Cells.Item(113,8).Value = 1.92
$ws.Cells.Item(113,9).Value = 1849300

# Row 114: 2019-11-26
$ws.Cells.Item(114,1).Value = 1574726400
$ws.Cells.Item(114,2).Value = "2019-11-26"
$ws.Cells.Item(114,3).Value = "0208"
$ws.Cells.Item(114,4).Value = "GREATEC"
$ws.Cells.Item(114,5).Value = 1.91
$ws.Cells.Item(114,6).Value = 1.95
$ws.Cells.Item(114,7).Value = 1.9
$ws.Cells.Item(114,8).Value = 1.91
$ws.Cells.Item(114,9).Value = 2163700

# Row 115: 2019-11-27
$ws.Cells.Item(115,1).Value = 1574812800
$ws.Cells.Item(115,2).Value = "2019-11-27"
$ws.Cells.Item(115,3).Value = "0208"
$ws.Cells.Item(115,4).Value = "GREATEC"
$ws.Cells.Item(115,5).Value = 1.91
$ws.Cells.Item(115,6).Value = 1.94
$ws.Cells.Item(115,7).Value = 1.83
$ws.Cells.Item(115,8).Value = 1.87
$ws.Cells.Item(115,9).Value = 3322300

# Row 116: 2019-11-28
$ws.Cells.Item(116,1).Value = 1574899200
$ws.Cells.Item(116,2).Value = "2019-11-28"
$ws.Cells.Item(116,3).Value = "0208"
$ws.Cells.Item(116,4).Value = "GREATEC"
$ws.Cells.Item(116,5).Value = 1.87
$ws.Cells.Item(116,6).Value = 1.92
$ws.Cells.Item(116,7).Value = 1.84
$ws.Cells.Item(116,8).Value = 1.92
$ws.Cells.Item(116,9).Value = 1322300

# Restore General number format on the text columns now that values are set
$textRng.ClearFormats()
